$d = $word.ActiveDocument

# Position at the very end of the document body (after the final "}" paragraph)
$r = $d.Content
$r.Collapse(0)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0"/>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Default"/>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t>background</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t>-</w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t>image</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t xml:space="preserve"> az</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t xml:space="preserve"> URL CSS</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t xml:space="preserve">- </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t>és nem a HTML-fájlhoz viszonyított</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
    <w:t>!!!</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Default"/>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Default"/>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Bariol" w:hAnsi="Bariol" w:cs="Bariol"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@

[void]$r.InsertXML($xml)
